$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 443
    3  = 445
    4  = 447
    5  = 450
    6  = 451
    7  = 454
    8  = 456
    9  = 458
    10 = 459
    11 = 461
    12 = 463
    13 = 465
    14 = 467
    15 = 24
    16 = 37
    17 = 100
    18 = 108
    19 = 182
    20 = 241
    21 = 288
    22 = 327
    23 = 357
    24 = 377
    25 = 396
    26 = 426
    27 = 511
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
